$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.75
$ws.Range("I2").Value = 2.15
$ws.Range("J2").Value = 4.5
$ws.Range("L2").Value = 2.88
$ws.Range("AA2").Value = 2.05
$ws.Range("AB2").Value = 1.7
$ws.Range("AM2").Value = 501
$ws.Range("AN2").Value = 6
$ws.Range("AO2").Value = 9
$ws.Range("AQ2").Value = 19

# Row 3
$ws.Range("G3").Value = 1.7
$ws.Range("I3").Value = 6.25
$ws.Range("AF3").Value = 12
$ws.Range("AO3").Value = 29
$ws.Range("AP3").Value = 21

# Row 4
$ws.Range("G4").Value = 1.44
$ws.Range("N4").Value = 8
$ws.Range("AF4").Value = 9
$ws.Range("AS4").Value = 81

# Row 5
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7

# Row 6
$ws.Range("G6").Value = 2.8
$ws.Range("H6").Value = 3.3
$ws.Range("I6").Value = 2.45
$ws.Range("J6").Value = 3.4
$ws.Range("L6").Value = 3.1
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 15
$ws.Range("AF6").Value = 29
$ws.Range("AG6").Value = 21
$ws.Range("AI6").Value = 11
$ws.Range("AN6").Value = 9.5
$ws.Range("AO6").Value = 13
$ws.Range("AP6").Value = 10
$ws.Range("AQ6").Value = 23
$ws.Range("AR6").Value = 19
$ws.Range("AS6").Value = 26

# Row 9
$ws.Range("AM9").Value = 400
